$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 512-513, shifting existing rows 512:613 down to 514:613
$ws.Range("A512:R513").Insert(-4121)

# --- Row 512 (new) ---
$ws.Cells.Item(512, 1).Value = 11
$ws.Cells.Item(512, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(512, 3).Value = "Bíobío"
$ws.Cells.Item(512, 4).Value = 45244
$ws.Cells.Item(512, 5).Value = 8
$ws.Cells.Item(512, 6).Value = 100112023
$ws.Cells.Item(512, 7).Value = "Brócoli"
$ws.Cells.Item(512, 8).Value = "Sin especificar"
$ws.Cells.Item(512, 9).Value = "Primera"
$ws.Cells.Item(512, 10).Value = 1000
$ws.Cells.Item(512, 11).Value = 1000
$ws.Cells.Item(512, 12).Value = 1000
$ws.Cells.Item(512, 13).Value = 1000
$ws.Cells.Item(512, 14).Value = "$/unidad"
$ws.Cells.Item(512, 15).Value = "Región Metropolitana"
$ws.Cells.Item(512, 16).Value = 1000
$ws.Cells.Item(512, 17).Value = 1
$ws.Cells.Item(512, 18).Value = "Hortaliza"

# --- Row 513 (new) ---
$ws.Cells.Item(513, 1).Value = 11
$ws.Cells.Item(513, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(513, 3).Value = "Bíobío"
$ws.Cells.Item(513, 4).Value = 45244
$ws.Cells.Item(513, 5).Value = 8
$ws.Cells.Item(513, 6).Value = 100112023
$ws.Cells.Item(513, 7).Value = "Brócoli"
$ws.Cells.Item(513, 8).Value = "Sin especificar"
$ws.Cells.Item(513, 9).Value = "Segunda"
$ws.Cells.Item(513, 10).Value = 1000
$ws.Cells.Item(513, 11).Value = 800
$ws.Cells.Item(513, 12).Value = 800
$ws.Cells.Item(513, 13).Value = 800
$ws.Cells.Item(513, 14).Value = "$/unidad"
$ws.Cells.Item(513, 15).Value = "Región Metropolitana"
$ws.Cells.Item(513, 16).Value = 800
$ws.Cells.Item(513, 17).Value = 1
$ws.Cells.Item(513, 18).Value = "Hortaliza"
